# Updated cryptos list on Tue Mar  7 16:34:39 UTC 2023 with GitHub Actions
#
# Updates the Price (D) and Volume(1h) (E) columns for each row, and for a
# handful of rows (two coins that swapped rank order) also updates the
# Coin (B) and Link (C) columns.
#
# Price values that look like plain numbers (e.g. "1.003", "5.830",
# "11.60") are written with a leading apostrophe so Excel stores them as
# literal text (preserving exact formatting/trailing zeros) instead of
# re-interpreting them as numbers; the cell Style is then reset to
# "Normal" so no stray quote-prefix formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.264.76'
$ws.Range('E2').Value = '  -1.22%  '
$ws.Range('D3').Value = '1.555.06'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''1.002'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').Value = '''286.77'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').Value = '''0.3772'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.63%  '
$ws.Range('D8').Value = '''0.3264'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.59%  '
$ws.Range('D9').Value = '''43.56'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -10.22%  '
$ws.Range('D10').Value = '''1.132'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D11').Value = '''0.07358'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = '''1.004'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('D13').Value = '''20.16'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -4.03%  '
$ws.Range('D14').Value = '''5.830'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.94%  '
$ws.Range('D15').Value = '''6.762'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('D16').Value = '1.558.22'
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('D17').Value = '''0.00001075'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -3.66%  '
$ws.Range('D18').Value = '''0.06638'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').Value = '''85.88'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '''1.002'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''6.369'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.83%  '
$ws.Range('D22').Value = '''16.08'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -2.76%  '
$ws.Range('D23').Value = '''11.60'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.67%  '
$ws.Range('D24').Value = '22.272.89'
$ws.Range('E24').Value = '  -1.16%  '
$ws.Range('D25').Value = '''2.300'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.34%  '
$ws.Range('D26').Value = '''2.554'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('D27').Value = '''150.14'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.74%  '
$ws.Range('D28').Value = '''19.33'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.70%  '
$ws.Range('D29').Value = '''4.927'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('D30').Value = '''122.29'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.64%  '
$ws.Range('D31').Value = '1.733.65'
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('D32').Value = '''1.079'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('D33').Value = '''5.891'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -4.92%  '
$ws.Range('D34').Value = '''1.901'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -5.19%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = '''0.08262'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.67%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '''9.329'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.41%  '
$ws.Range('D37').Value = '''0.02365'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.75%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '''0.06295'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '''5.275'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.50%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').Value = '''0.2148'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -5.47%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '''1.254'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D42').Value = '''11.01'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').Value = '''1.001'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '''0.6039'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.77%  '
$ws.Range('D45').Value = '''13.67'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.31%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '''3.737'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.5889'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.56%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''123.59'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '''1.968'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.72%  '
$ws.Range('D50').Value = '''1.176'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.76%  '
$ws.Range('D51').Value = '''0.07076'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.81%  '
